# Window handle sample implementation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Nationality values for row 2 and row 3
$ws.Range("D2").Value = "America"
$ws.Range("D3").Value = "Srilanka"

# Update the selected cell to D4
$ws.Range("D4").Select()
